# Apply changes to the "BatchModule" sheet (4th sheet) as described by the diff:
#   C5: 123 (number)         -> "AWS1" (shared string, same style as C2:C4)
#   D5: 23 (number)          -> TRUE (boolean, same as D2:D4)
#   E5: "dats" (shared str)  -> 5 (number), keeping existing style

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BatchModule")

# Mirror the style used by the row above (row 4) for C5 (Arial 10pt black, same as C2:C4).
$ws.Range("C5").Value = "AWS1"
$ws.Range("C5").Font.Name = "Arial"
$ws.Range("C5").Font.Size = 10
$ws.Range("C5").Font.Color = 0

$ws.Range("D5").Value = $true

$ws.Range("E5").Value = 5

# Shrink the window height, matching the workbookView change in the diff.
$excel.ActiveWindow.Height = 21100
